# RequirementsRandomization.xlsx update:
#  - Remove the "Top.A.Y" sheet
#  - Rename "Top.B" -> "CRV.Constraints" and populate it
#  - Rename "Top.A.X" -> "CRV.General" and populate it
#  - Update CRV / CRV.DataTypes / CRV.Phases sheets with new/changed
#    requirement rows reflecting the DVCon US 2021 workshop feedback

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet bookkeeping: drop Top.A.Y, repurpose Top.B / Top.A.X
# ---------------------------------------------------------------------
$wsTopAY = $wb.Worksheets.Item("Top.A.Y")
$wsTopAY.Delete()

$wsConstraints = $wb.Worksheets.Item("Top.B")
$wsConstraints.Name = "CRV.Constraints"

$wsGeneral = $wb.Worksheets.Item("Top.A.X")
$wsGeneral.Name = "CRV.General"

# ---------------------------------------------------------------------
# 2. CRV sheet: update trace-down columns + add Constraints/General rows
# ---------------------------------------------------------------------
$wsCRV = $wb.Worksheets.Item("CRV")

# Row 28 (Class) / Row 29 (DataTypes): "trace down" now points at the
# dedicated requirement sheet instead of a comma separated tag list.
$wsCRV.Range("F28").Value = "CRV.Class"
$wsCRV.Range("F29").Value = "CRV.Datatypes"

# Row 30 (Phases): add the "trace down" pointer to its own sheet
$wsCRV.Range("F30").Value = "CRV.Phases"
$wsCRV.Range("K30").Value = "COSEDA"

# Row 31: new "Constraints" requirement group
$wsCRV.Range("A31").Value = "Constraints"
$wsCRV.Range("B31").Value = "Support randomization using constraints"
$wsCRV.Range("C31").Value = "Support for constraints"
$wsCRV.Range("D31").Value = "All kinds of randomization"
$wsCRV.Range("F31").Value = "CRV.Constraints"
$wsCRV.Range("K31").Value = "COSEDA, DVCON US 2021 Tutorial Remark "

# Row 32: new "General" requirement group
$wsCRV.Range("A32").Value = "General"
$wsCRV.Range("B32").Value = "General Requirements which are not part of other sections"
$wsCRV.Range("F32").Value = "CRV.General"
$wsCRV.Range("K32").Value = "COSEDA, DVCON US 2021 Tutorial Remark "

# ---------------------------------------------------------------------
# 3. CRV.DataTypes sheet: update DataTypes.2 + add DataTypes.6 row
# ---------------------------------------------------------------------
$wsDT = $wb.Worksheets.Item("CRV.DataTypes")

$wsDT.Range("C4").Value = "Real values are used AMS (SystemC AMS)modelling "
$wsDT.Range("D4").Value = "Analog / Real value verification" + [char]10 + "Also in Discusstion for UVM AMS on SystemVerilog "

$wsDT.Range("A8").Value = "DataTypes.6"
$wsDT.Range("B8").Value = "Support fixed point datatypes "
$wsDT.Range("C8").Value = "Extension to fixed point demantic"
$wsDT.Range("D8").Value = "UVM standard payload randomization"
$wsDT.Range("K8").Value = "DVCON US 2021 Tutorial Remark "

# ---------------------------------------------------------------------
# 4. CRV.Phases sheet: add Phases.2 row
# ---------------------------------------------------------------------
$wsPhases = $wb.Worksheets.Item("CRV.Phases")

$wsPhases.Range("A4").Value = "Phases.2"
$wsPhases.Range("B4").Value = "Support pre_randomize() and post_randomize() callback"
$wsPhases.Range("C4").Value = "Enables callback for pre and post randomization"
$wsPhases.Range("K4").Value = "COSEDA"

# ---------------------------------------------------------------------
# 5. CRV.Constraints sheet (was Top.B): header row + two requirements
# ---------------------------------------------------------------------
$wsConstraints.Range("A2").Value = "Tag"
$wsConstraints.Range("B2").Value = "Requirement"
$wsConstraints.Range("C2").Value = "Rationale"
$wsConstraints.Range("D2").Value = "Use Cases"
$wsConstraints.Range("E2").Value = "Trace Up"
$wsConstraints.Range("F2").Value = "Trace Down"
$wsConstraints.Range("G2").Value = "Stability (1 least - 10 most)"
$wsConstraints.Range("H2").Value = "Sope (H, M, L)"
$wsConstraints.Range("I2").Value = "Priority (H, M, L)"
$wsConstraints.Range("J2").Value = "Type (P, F, E, Q)"
$wsConstraints.Range("K2").Value = "Source"

$wsConstraints.Range("A3").Value = "Constraints.1"
$wsConstraints.Range("B3").Value = "Support for soft constraints"
$wsConstraints.Range("C3").Value = "Enables the ability to add  constraints whcih can be dropped to avoid overconstrainging"
$wsConstraints.Range("D3").Value = "UVM standard payload randomization"
$wsConstraints.Range("K3").Value = "DVCON US 2021 Tutorial Remark "

$wsConstraints.Range("A4").Value = "Constraints.2"
$wsConstraints.Range("B4").Value = "Ability to Debug conflicting constraints"
$wsConstraints.Range("C4").Value = "Conflicting constraints can easily occur"
$wsConstraints.Range("D4").Value = "More complex constraints"
$wsConstraints.Range("K4").Value = "DVCON US 2021 Tutorial Remark "

# Requirement-column header comments (same boilerplate used on the other
# requirement sheets)
$wsConstraints.Range("B2").AddComment("Succinct requirement decsription") | Out-Null
$wsConstraints.Range("C2").AddComment("Why requirement exists") | Out-Null
$wsConstraints.Range("E2").AddComment("higher level requirement this requirement is considered to be a sub-requirement of") | Out-Null
$wsConstraints.Range("F2").AddComment("Next level sub-requirement that this requirement is considered to be a parent or the architectural component(s) that this last sub level requirement maps to") | Out-Null
$wsConstraints.Range("G2").AddComment("Likelihood requirement will not change") | Out-Null
$wsConstraints.Range("H2").AddComment("Extent to which this requirement interacts with other requirements or results in architectural changes") | Out-Null
$wsConstraints.Range("I2").AddComment("Importance of meeting this requirement in order relative to other requirements") | Out-Null
$wsConstraints.Range("J2").AddComment("Nature of requirement as a functional requirement vs non-functional, quantifiable vs non-quantifiable, emerging vs non-emerging, process vs non-emerging") | Out-Null
$wsConstraints.Range("K2").AddComment("Stakeholders that contributed the requirement") | Out-Null

# ---------------------------------------------------------------------
# 6. CRV.General sheet (was Top.A.X): header row + one requirement
# ---------------------------------------------------------------------
$wsGeneral.Range("A2").Value = "Tag"
$wsGeneral.Range("B2").Value = "Requirement"
$wsGeneral.Range("C2").Value = "Rationale"
$wsGeneral.Range("D2").Value = "Use Cases"
$wsGeneral.Range("E2").Value = "Trace Up"
$wsGeneral.Range("F2").Value = "Trace Down"
$wsGeneral.Range("G2").Value = "Stability (1 least - 10 most)"
$wsGeneral.Range("H2").Value = "Sope (H, M, L)"
$wsGeneral.Range("I2").Value = "Priority (H, M, L)"
$wsGeneral.Range("J2").Value = "Type (P, F, E, Q)"
$wsGeneral.Range("K2").Value = "Source"

$wsGeneral.Range("A3").Value = "General.1"
$wsGeneral.Range("B3").Value = "Support for randomization with reproducable results (e.g. based on a seed)"
$wsGeneral.Range("C3").Value = "Enables the ability to add  constraints whcih can be dropped to avoid overconstrainging"
$wsGeneral.Range("D3").Value = "UVM standard payload randomization"
$wsGeneral.Range("K3").Value = "DVCON US 2021 Tutorial Remark "

$wsGeneral.Range("B2").AddComment("Succinct requirement decsription") | Out-Null
$wsGeneral.Range("C2").AddComment("Why requirement exists") | Out-Null
$wsGeneral.Range("E2").AddComment("higher level requirement this requirement is considered to be a sub-requirement of") | Out-Null
$wsGeneral.Range("F2").AddComment("Next level sub-requirement that this requirement is considered to be a parent or the architectural component(s) that this last sub level requirement maps to") | Out-Null
$wsGeneral.Range("G2").AddComment("Likelihood requirement will not change") | Out-Null
$wsGeneral.Range("H2").AddComment("Extent to which this requirement interacts with other requirements or results in architectural changes") | Out-Null
$wsGeneral.Range("I2").AddComment("Importance of meeting this requirement in order relative to other requirements") | Out-Null
$wsGeneral.Range("J2").AddComment("Nature of requirement as a functional requirement vs non-functional, quantifiable vs non-quantifiable, emerging vs non-emerging, process vs non-emerging") | Out-Null
$wsGeneral.Range("K2").AddComment("Stakeholders that contributed the requirement") | Out-Null

# ---------------------------------------------------------------------
# 7. Re-select CRV as the active/visible sheet (matches tabSelected)
# ---------------------------------------------------------------------
$wsCRV.Activate()
$wsCRV.Range("K32").Select() | Out-Null
